# importTenTen.xlsx — apply the "chang excel then builded" edit:
#  - C4: "Doãn Trí Bình" -> "SAM SAM"
#  - F3: "lon.com"       -> "skybabies.com"
#  - F6: "lch.com"       -> "lichntb.com"
#  - New column H: header "郵便番号" (H1) and 10000 in H2:H6
#  - Size column H and move the selection to H6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content fixes on existing columns (order matters for shared-string layout) ---
$ws.Range("C4").Value = "SAM SAM"
$ws.Range("F3").Value = "skybabies.com"
$ws.Range("F6").Value = "lichntb.com"

# --- new column H ("郵便番号") ---
# Copy the header formatting from G1 (bold/fill/centered header style) onto H1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "郵便番号"
$ws.Range("H1").ColumnWidth = 20.5

$ws.Range("H2").Value = 10000
$ws.Range("H3").Value = 10000
$ws.Range("H4").Value = 10000
$ws.Range("H5").Value = 10000
$ws.Range("H6").Value = 10000

# --- final selection as left in the saved file ---
$ws.Range("H6").Select()
